$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.216.05"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.429.32"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.432.20"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  -8.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").Value = "4.021.20"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("E16").Value = "  -6.62%  "
$ws.Range("D17").Value = "64.292.47"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "3.423.44"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000116"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.859"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.43%  "
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.814.34"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.65%  "
$ws.Range("E48").Value = "  +8.18%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  -4.09%  "
